$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.253.66'
$ws.Range("E2").Value = '  +5.84%  '

$ws.Range("D3").Value = '2.388.99'
$ws.Range("E3").Value = '  +4.38%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''550.95'
$ws.Range("E5").Value = '  +2.97%  '

$ws.Range("D6").Value = '''134.85'
$ws.Range("E6").Value = '  +2.56%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +3.42%  '

$ws.Range("D9").Value = '2.385.90'
$ws.Range("E9").Value = '  +4.22%  '

$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("D11").Value = '''5.57'
$ws.Range("E11").Value = '  +2.36%  '

$ws.Range("E12").Value = '  +1.82%  '

$ws.Range("E13").Value = '  +3.05%  '

$ws.Range("D14").Value = '''24.44'
$ws.Range("E14").Value = '  +4.16%  '

$ws.Range("D15").Value = '2.810.80'
$ws.Range("E15").Value = '  +4.29%  '

$ws.Range("D16").Value = '61.096.22'
$ws.Range("E16").Value = '  +5.59%  '

$ws.Range("E17").Value = '  +2.48%  '

$ws.Range("D18").Value = '2.285.42'
$ws.Range("E18").Value = '  +0.57%  '

$ws.Range("D19").Value = '''10.86'
$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''7.07'
$ws.Range("E20").Value = '  +10.71%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '''4.25'
$ws.Range("E21").Value = '  +0.42%  '

$ws.Range("D22").Value = '''320.89'
$ws.Range("E22").Value = '  +2.59%  '

$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").Value = '''63.92'
$ws.Range("E24").Value = '  +1.78%  '

$ws.Range("E25").Value = '  +5.16%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("D27").Value = '''8.22'
$ws.Range("E27").Value = '  +3.32%  '

$ws.Range("E28").Value = '  +4.11%  '

$ws.Range("E29").Value = '  +2.60%  '

$ws.Range("E30").Value = '  +4.89%  '

$ws.Range("D31").Value = '''172.10'
$ws.Range("E31").Value = '  +1.28%  '

$ws.Range("E32").Value = '  +7.73%  '

$ws.Range("D33").Value = '''6.00'
$ws.Range("E33").Value = '  +4.57%  '

$ws.Range("E34").Value = '  +15.68%  '

$ws.Range("D35").Value = '''0.389'
$ws.Range("E35").Value = '  +2.72%  '

$ws.Range("D36").Value = '''18.22'
$ws.Range("E36").Value = '  +3.23%  '

$ws.Range("D38").Value = '''4.27'
$ws.Range("E38").Value = '  +9.68%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").Value = '''325.86'
$ws.Range("E40").Value = '  +12.39%  '

$ws.Range("E41").Value = '  +6.85%  '

$ws.Range("E42").Value = '  +1.23%  '

$ws.Range("D43").Value = '''146.11'
$ws.Range("E43").Value = '  +2.46%  '

$ws.Range("D44").Value = '''3.52'
$ws.Range("E44").Value = '  +3.56%  '

$ws.Range("E45").Value = '  +1.48%  '

$ws.Range("E46").Value = '  +8.80%  '

$ws.Range("D47").Value = '''0.0508'
$ws.Range("E47").Value = '  +2.27%  '

$ws.Range("E48").Value = '  +2.61%  '

$ws.Range("D49").Value = '''0.0216'
$ws.Range("E49").Value = '  +2.99%  '

$ws.Range("D50").Value = '''11.03'

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0201'
$ws.Range("E51").Value = '  +2.22%  '
